# Auto update Excel log 2026-02-04 14:14:47
# Appends newly captured sensor readings (2026-02-04, ~14:13-14:14) to the
# PIR, Humidity, and Temperature logs.
#
# Each data row is flattened into 7 consecutive array slots in the order:
#   RowNumber, Date, Timestamp, Hour, Location, Value, Status
# (this interpreter does not reliably support nested/jagged array literals,
# so a flat array with a fixed stride of 7 is used instead).

$wb = $excel.ActiveWorkbook
$stride = 7

# --- PIR sheet: add rows 153-164 (dimension A1:F152 -> A1:F164) ---
$pirData = @(
    "153", "2026-02-04", "14:13:43", "14:00", "Bathroom", "No Motion", "Inactive",
    "154", "2026-02-04", "14:13:44", "14:00", "Bathroom", "No Motion", "Inactive",
    "155", "2026-02-04", "14:13:49", "14:00", "Bathroom", "No Motion", "Inactive",
    "156", "2026-02-04", "14:13:54", "14:00", "Bathroom", "No Motion", "Inactive",
    "157", "2026-02-04", "14:13:59", "14:00", "Bathroom", "No Motion", "Inactive",
    "158", "2026-02-04", "14:14:04", "14:00", "Bathroom", "No Motion", "Inactive",
    "159", "2026-02-04", "14:14:09", "14:00", "Bathroom", "No Motion", "Inactive",
    "160", "2026-02-04", "14:14:14", "14:00", "Bathroom", "No Motion", "Inactive",
    "161", "2026-02-04", "14:14:19", "14:00", "Bathroom", "No Motion", "Inactive",
    "162", "2026-02-04", "14:14:24", "14:00", "Bathroom", "Motion Detected", "Active",
    "163", "2026-02-04", "14:14:32", "14:00", "Bathroom", "No Motion", "Inactive",
    "164", "2026-02-04", "14:14:37", "14:00", "Bathroom", "No Motion", "Inactive"
)
$ws = $wb.Worksheets.Item("PIR")
$count = $pirData.Count / $stride
for ($i = 0; $i -lt $count; $i++) {
    $base = $i * $stride
    $rowNum = $pirData[$base]

    # Column A holds date strings like "2026-02-04"; force text format so
    # Excel does not auto-convert it into a date serial number.
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $pirData[$base + 1]
    $ws.Cells.Item($rowNum, 2).Value = $pirData[$base + 2]
    $ws.Cells.Item($rowNum, 3).Value = $pirData[$base + 3]
    $ws.Cells.Item($rowNum, 4).Value = $pirData[$base + 4]
    $ws.Cells.Item($rowNum, 5).Value = $pirData[$base + 5]
    $ws.Cells.Item($rowNum, 6).Value = $pirData[$base + 6]
}

# --- Humidity sheet: add rows 122-135 (dimension A1:F121 -> A1:F135) ---
$humidityData = @(
    "122", "2026-02-04", "14:13:42", "14:00", "Bathroom", "76.4%", "Active",
    "123", "2026-02-04", "14:13:43", "14:00", "Bathroom", "77.4%", "Active",
    "124", "2026-02-04", "14:13:46", "14:00", "Bathroom", "76.3%", "Active",
    "125", "2026-02-04", "14:13:51", "14:00", "Bathroom", "77.2%", "Active",
    "126", "2026-02-04", "14:13:56", "14:00", "Bathroom", "76.2%", "Active",
    "127", "2026-02-04", "14:14:01", "14:00", "Bathroom", "77.2%", "Active",
    "128", "2026-02-04", "14:14:06", "14:00", "Bathroom", "76.3%", "Active",
    "129", "2026-02-04", "14:14:11", "14:00", "Bathroom", "77.3%", "Active",
    "130", "2026-02-04", "14:14:16", "14:00", "Bathroom", "76.4%", "Active",
    "131", "2026-02-04", "14:14:21", "14:00", "Bathroom", "77.3%", "Active",
    "132", "2026-02-04", "14:14:26", "14:00", "Bathroom", "76.6%", "Active",
    "133", "2026-02-04", "14:14:31", "14:00", "Bathroom", "77.5%", "Active",
    "134", "2026-02-04", "14:14:36", "14:00", "Bathroom", "76.7%", "Active",
    "135", "2026-02-04", "14:14:41", "14:00", "Bathroom", "77.7%", "Active"
)
$ws = $wb.Worksheets.Item("Humidity")
$count = $humidityData.Count / $stride
for ($i = 0; $i -lt $count; $i++) {
    $base = $i * $stride
    $rowNum = $humidityData[$base]

    # Column A holds date strings like "2026-02-04"; force text format so
    # Excel does not auto-convert it into a date serial number.
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $humidityData[$base + 1]
    $ws.Cells.Item($rowNum, 2).Value = $humidityData[$base + 2]
    $ws.Cells.Item($rowNum, 3).Value = $humidityData[$base + 3]
    $ws.Cells.Item($rowNum, 4).Value = $humidityData[$base + 4]

    # Value column holds humidity strings like "76.4%"; force text so it
    # does not get auto-converted into a percentage number.
    $ws.Cells.Item($rowNum, 5).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 5).Value = $humidityData[$base + 5]
    $ws.Cells.Item($rowNum, 6).Value = $humidityData[$base + 6]
}

# --- Temperature sheet: add rows 122-135 (dimension A1:F121 -> A1:F135) ---
$temperatureData = @(
    "122", "2026-02-04", "14:13:42", "14:00", "Bathroom", "24.9C", "Active",
    "123", "2026-02-04", "14:13:44", "14:00", "Bathroom", "24.9C", "Active",
    "124", "2026-02-04", "14:13:46", "14:00", "Bathroom", "24.9C", "Active",
    "125", "2026-02-04", "14:13:51", "14:00", "Bathroom", "24.9C", "Active",
    "126", "2026-02-04", "14:13:56", "14:00", "Bathroom", "24.9C", "Active",
    "127", "2026-02-04", "14:14:02", "14:00", "Bathroom", "24.9C", "Active",
    "128", "2026-02-04", "14:14:07", "14:00", "Bathroom", "24.9C", "Active",
    "129", "2026-02-04", "14:14:12", "14:00", "Bathroom", "24.9C", "Active",
    "130", "2026-02-04", "14:14:17", "14:00", "Bathroom", "24.8C", "Active",
    "131", "2026-02-04", "14:14:22", "14:00", "Bathroom", "24.8C", "Active",
    "132", "2026-02-04", "14:14:27", "14:00", "Bathroom", "24.8C", "Active",
    "133", "2026-02-04", "14:14:32", "14:00", "Bathroom", "24.8C", "Active",
    "134", "2026-02-04", "14:14:37", "14:00", "Bathroom", "24.8C", "Active",
    "135", "2026-02-04", "14:14:42", "14:00", "Bathroom", "24.8C", "Active"
)
$ws = $wb.Worksheets.Item("Temperature")
$count = $temperatureData.Count / $stride
for ($i = 0; $i -lt $count; $i++) {
    $base = $i * $stride
    $rowNum = $temperatureData[$base]

    # Column A holds date strings like "2026-02-04"; force text format so
    # Excel does not auto-convert it into a date serial number.
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $temperatureData[$base + 1]
    $ws.Cells.Item($rowNum, 2).Value = $temperatureData[$base + 2]
    $ws.Cells.Item($rowNum, 3).Value = $temperatureData[$base + 3]
    $ws.Cells.Item($rowNum, 4).Value = $temperatureData[$base + 4]
    $ws.Cells.Item($rowNum, 5).Value = $temperatureData[$base + 5]
    $ws.Cells.Item($rowNum, 6).Value = $temperatureData[$base + 6]
}

